$wb = $excel.ActiveWorkbook

# "展览" sheet (sheet1) and "全部类型" sheet (sheet4) both carry the same
# event data; update the "想去人数" (want-to-go count) column F for rows 2 and 3.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 690
    $ws.Range("F3").Value = 4004
}
